# Apply the "456a3b4" gh-pages data refresh to 杭州-漫展信息.xlsx
# Sheets (tab order): 1=展览, 2=演出, 3=本地生活, 4=全部类型

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # 展览
$ws2 = $wb.Worksheets.Item(2)   # 演出
$ws3 = $wb.Worksheets.Item(3)   # 本地生活
$ws4 = $wb.Worksheets.Item(4)   # 全部类型

# ---------------------------------------------------------------
# Sheet 1 (展览): refreshed "想去人数" (F) counts + one cover URL
# ---------------------------------------------------------------
$ws1.Range("F6").Value = 43
$ws1.Range("F10").Value = 10227
$ws1.Range("F11").Value = 182
$ws1.Range("F12").Value = 71
$ws1.Range("F15").Value = 1963
$ws1.Range("F20").Value = 171
$ws1.Range("F23").Value = 1116
$ws1.Range("F24").Value = 83
$ws1.Range("F26").Value = 626
$ws1.Range("F28").Value = 154
$ws1.Range("I28").Value = "//i2.hdslb.com/bfs/openplatform/202409/5choDLVP1726713753891.png"
$ws1.Range("F29").Value = 626
$ws1.Range("F30").Value = 2773
$ws1.Range("F31").Value = 947
$ws1.Range("F32").Value = 681
$ws1.Range("F36").Value = 530
$ws1.Range("F37").Value = 207
$ws1.Range("F38").Value = 12
$ws1.Range("F39").Value = 1233
$ws1.Range("F40").Value = 558
$ws1.Range("F41").Value = 5246
$ws1.Range("F43").Value = 68
$ws1.Range("F44").Value = 108
$ws1.Range("F45").Value = 165
$ws1.Range("F47").Value = 4048

# ---------------------------------------------------------------
# Sheet 2 (演出): a new concert was inserted as row 13, pushing the
# previously-tracked rows 13-22 down to rows 14-23 (dimension grows
# from A1:I22 to A1:I23). Insert a row and populate it.
# ---------------------------------------------------------------
$ws2.Rows.Item(13).Insert()

# Row 13 lost its border/bold formatting on insert - pull it back from
# the row that now sits right below (the old row 13).
$ws2.Range("A14:I14").Copy()
$ws2.Range("A13:I13").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = 0

$ws2.Range("A13").Value = 12
# Force column B to stay plain text (matches every other row) instead
# of being auto-parsed into a date serial.
$ws2.Range("B13").NumberFormat = "@"
$ws2.Range("B13").Value = "2024-11-08"
$ws2.Range("C13").Value = "杭州·『YOKO TAKAHASHI EVANGELION ultimate Live「月十夜」』EVA高桥洋子巡演"
$ws2.Range("D13").Value = "武林之星博览中心2号楼 杭州 unilivehouse  （由你现场）"
$ws2.Range("E13").Value = "2024.11.08 20:00-11.08 21:30"
$ws2.Range("F13").Value = 20
$ws2.Range("G13").Value = 280
$ws2.Range("H13").Value = "https://show.bilibili.com/platform/detail.html?id=92580"
$ws2.Range("I13").Value = "//i1.hdslb.com/bfs/openplatform/202409/5bGHjiLT1726653844731.jpeg"

# ---------------------------------------------------------------
# Sheet 3 (本地生活): refreshed "想去人数" for the first listing
# ---------------------------------------------------------------
$ws3.Range("F2").Value = 732

# ---------------------------------------------------------------
# Sheet 4 (全部类型): same refreshed counts as sheets 1-3, combined
# ---------------------------------------------------------------
$ws4.Range("F2").Value = 732
$ws4.Range("F9").Value = 43
$ws4.Range("F13").Value = 10227
$ws4.Range("F14").Value = 182
$ws4.Range("F15").Value = 71
$ws4.Range("F22").Value = 1116
$ws4.Range("F23").Value = 83
$ws4.Range("F27").Value = 626
$ws4.Range("F29").Value = 154
$ws4.Range("I29").Value = "//i2.hdslb.com/bfs/openplatform/202409/5choDLVP1726713753891.png"
$ws4.Range("F30").Value = 626
$ws4.Range("F31").Value = 2774
$ws4.Range("F32").Value = 947
$ws4.Range("F37").Value = 207
$ws4.Range("F39").Value = 68
$ws4.Range("F40").Value = 108
$ws4.Range("F41").Value = 165
$ws4.Range("F43").Value = 4048
